# SBA Science Variable Dictionary 2018-2019.xlsx
#
# 1. Rename the first sheet from "SBAScience2019" to "SBASciSPRING1819".
# 2. On "SBASciFALL1819", remove the obsolete "testname" row (old row 13)
#    and three other rows that no longer apply (old rows 20-22), letting
#    every row below shift up to close the gap.
# 3. Restore the active-cell selections left in each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "SBASciSPRING1819"

$ws2 = $wb.Worksheets.Item(2)

# Delete from the bottom up so the remaining row numbers aren't
# invalidated by earlier deletions.
$ws2.Rows.Item(22).Delete()
$ws2.Rows.Item(21).Delete()
$ws2.Rows.Item(20).Delete()
$ws2.Rows.Item(13).Delete()

# Match the saved cursor/selection position on each sheet.
[void]$ws1.Range("E20").Select()
[void]$ws2.Range("E27").Select()
